$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 'Did you know that cats can make over 100 different sounds? Dogs, on the other hand, can barely muster a couple of "woofs."'
    3  = "Fact: Penguins are great comedians. Isn't it amazing how they can be so funny without even cracking a smile?"
    4  = "Did you know that dolphins are always smiling? It's probably because they know they're the smartest fish in the sea!"
    5  = "Fact: The average American eats around 35 tons of food in their lifetime. But my guess is that half of that is just pizza and ice cream!"
    6  = "Did you know that frogs are great jumpers? They're basically the Olympic athletes of the amphibian world!"
    7  = "Fact: The average person laughs around 13 times a day. Unless, of course, they forget to wear their funny socks—then it drops down to zero!"
    8  = "Did you know that squirrels forget where they bury about half of their nuts? That's why we find them in the most unexpected places!"
    9  = "Fact: It takes an average of 364 licks to get to the center of a Tootsie Pop. But who has the patience for that? Just take a bite!"
    10 = "Did you know that ducks have feathers to cover their butt quacks? It's true! That way, they can keep their tails wagging with dignity!"
    11 = "Fact: The world's largest recorded snowflake measured 15 inches in diameter. That must have been a very confused snowman!"
    12 = "Did you know that mosquitos are the deadliest animals in the world? They're just trying to tell us to invest in bug spray!"
    13 = "Fact: The average person spends around six months of their lifetime waiting at red lights. That's a lot of car karaoke time wasted!"
    14 = "Did you know that honey never spoils? So the next time someone makes a bee pun, just tell them that joke will never get old!"
    15 = "Fact: Banging your head against a wall for one hour burns 150 calories. But I'm not sure it's the most effective workout routine!"
    16 = "Did you know that elephants are the only animals that can't jump? No wonder they're always so grounded and accommodating!"
    17 = "Fact: The average person blinks around 17,000 times a day. That's a lot of missed opportunities to admire cute animals and epic sunsets!"
    18 = "Did you know that the first oranges weren't actually orange? They were green like limes! Talk about a citrus transformation!"
    19 = 'Fact: The slogan "Imagination at Work" was once considered by deodorant brand Axe. Can you imagine? Smelling like a wild adventure!'
    20 = "Did you know that cows have best friends? They udderly cherish their companionship, making everyday moosic together!"
    21 = "Fact: The world's largest recorded pizza measured 131 feet in diameter. Imagine the size of the delivery guy's car!"
}

foreach ($row in $values.Keys) {
    $ws.Range("A$row").Value = $values[$row]
}
